# Add working set of sequences
# For every data row where column E is TRUE and column F already contains
# "N/A", fill columns G through N (image1, type1, image2, type2, image3,
# type3, image4, type4) with "N/A" as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $eVal = $ws.Cells.Item($r, 5).Value()
    $fVal = $ws.Cells.Item($r, 6).Value()
    if ($eVal -eq $true -and $fVal -eq "N/A") {
        for ($c = 7; $c -le 14; $c++) {
            $ws.Cells.Item($r, $c).Value = "N/A"
        }
    }
}
